# Finished Week 13 logging
# Update row 3 ("R" row) values on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 128
$wsOff.Range("C3").Value = 87
$wsOff.Range("D3").Value = 35
$wsOff.Range("E3").Value = 21
$wsOff.Range("F3").Value = 1
$wsOff.Range("G3").Value = 2

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 116
$wsDef.Range("C3").Value = 94
$wsDef.Range("D3").Value = 27
$wsDef.Range("E3").Value = 10
$wsDef.Range("F3").Value = 1
$wsDef.Range("G3").Value = 1
